$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.824.95'
$ws.Range("E2").Value = '  +1.80%  '

$ws.Range("D3").Value = '3.331.83'
$ws.Range("E3").Value = '  +1.99%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = '581.62'
$ws.Range("E5").Value = '  +1.17%  '

$ws.Range("D6").Value = '176.85'
$ws.Range("E6").Value = '  +3.04%  '

$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").Value = '0.589'
$ws.Range("E8").Value = '  +2.46%  '

$ws.Range("D9").Value = '3.329.31'
$ws.Range("E9").Value = '  +2.18%  '

$ws.Range("D10").Value = '0.184'
$ws.Range("E10").Value = '  +7.53%  '

$ws.Range("D11").Value = '0.582'
$ws.Range("E11").Value = '  +2.67%  '

$ws.Range("D12").Value = '47.15'
$ws.Range("E12").Value = '  +5.45%  '

$ws.Range("D13").Value = '0.0000273'
$ws.Range("E13").Value = '  +3.41%  '

$ws.Range("D14").Value = '699.75'
$ws.Range("E14").Value = '  +1.01%  '

$ws.Range("D15").Value = '3.884.73'
$ws.Range("E15").Value = '  +2.52%  '

$ws.Range("D16").Value = '8.42'
$ws.Range("E16").Value = '  +2.68%  '

$ws.Range("D17").Value = '67.893.95'
$ws.Range("E17").Value = '  +1.81%  '

$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("D19").Value = '3.348.89'
$ws.Range("E19").Value = '  +2.61%  '

$ws.Range("D20").Value = '17.48'
$ws.Range("E20").Value = '  +2.09%  '

$ws.Range("D21").Value = '11.05'
$ws.Range("E21").Value = '  +4.24%  '

$ws.Range("D22").Value = '0.894'
$ws.Range("E22").Value = '  +1.90%  '

$ws.Range("E23").Value = '  +4.65%  '

$ws.Range("D24").Value = '17.02'
$ws.Range("E24").Value = '  +1.25%  '

$ws.Range("D25").Value = '99.99'
$ws.Range("E25").Value = '  +0.37%  '

$ws.Range("D26").Value = '3.92'
$ws.Range("E26").Value = '  +2.47%  '

$ws.Range("D27").Value = '2.70'
$ws.Range("E27").Value = '  +1.91%  '

$ws.Range("D28").Value = '9.56'
$ws.Range("E28").Value = '  +6.04%  '

$ws.Range("D29").Value = '33.01'
$ws.Range("E29").Value = '  -0.94%  '

$ws.Range("D30").Value = '8.57'
$ws.Range("E30").Value = '  +3.61%  '

$ws.Range("D31").Value = '7.01'
$ws.Range("E31").Value = '  +6.22%  '

$ws.Range("D32").Value = '567.78'
$ws.Range("E32").Value = '  -0.73%  '

$ws.Range("D33").Value = '10.99'
$ws.Range("E33").Value = '  +2.41%  '

$ws.Range("D34").Value = '0.105'
$ws.Range("E34").Value = '  +3.77%  '

$ws.Range("D35").Value = '57.53'
$ws.Range("E35").Value = '  +4.63%  '

$ws.Range("E36").Value = '  -0.03%  '

$ws.Range("D37").Value = '3.684.85'
$ws.Range("E37").Value = '  -3.74%  '

$ws.Range("D38").Value = '3.38'
$ws.Range("E38").Value = '  +2.26%  '

$ws.Range("D39").Value = '34.71'
$ws.Range("E39").Value = '  +11.35%  '

$ws.Range("E40").Value = '  +4.89%  '

$ws.Range("D41").Value = '2.64'
$ws.Range("E41").Value = '  +3.35%  '

$ws.Range("D42").Value = '3.17'
$ws.Range("E42").Value = '  +7.48%  '

$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").Value = '0.0₃0675'
$ws.Range("E43").Value = '  +2.83%  '

$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").Value = '0.336'
$ws.Range("E44").Value = '  +4.26%  '

$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").Value = '3.32'
$ws.Range("E45").Value = '  -1.07%  '

$ws.Range("D46").Value = '0.0408'
$ws.Range("E46").Value = '  +2.52%  '

$ws.Range("D47").Value = '2.67'
$ws.Range("E47").Value = '  +6.31%  '

$ws.Range("E48").Value = '  +1.88%  '

$ws.Range("E49").Value = '  +0.01%  '

$ws.Range("E50").Value = '  -0.38%  '

$ws.Range("D51").Value = '130.89'
$ws.Range("E51").Value = '  +1.28%  '
